# Split the single "employment" sheet into three sheets:
#   employment_smales   (unchanged 0.6 values)
#   employment_sfemales (values updated to 0.4 - "working activity status alignment")
#   employment_couples  (unchanged 0.6 values, copied from the original)
# and make employment_sfemales the active sheet, matching the target workbook layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("employment")

# Create the "sfemales" copy right after the original sheet.
$ws.Copy($null, $ws)
$sfemales = $wb.Worksheets.Item("employment (2)")

# Create the "couples" copy right after the "sfemales" copy.
$ws.Copy($null, $sfemales)
$couples = $wb.Worksheets.Item("employment (3)")

# Rename all three sheets to their final names.
$ws.Name = "employment_smales"
$sfemales.Name = "employment_sfemales"
$couples.Name = "employment_couples"

# Update the employed_share values (column B, rows 2-19) on the "sfemales" sheet to 0.4.
for ($r = 2; $r -le 19; $r++) {
    $sfemales.Cells.Item($r, 2).Value = 0.4
}

# Make "employment_sfemales" the active/selected sheet (matches activeTab + tabSelected in target).
$sfemales.Activate()
$sfemales.Range("B2").Select()
